$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting rows 12-17 down to 13-18
$ws.Rows.Item(12).Insert()

# Fill the new row 12 with the "Recurring meeting" entry
$ws.Cells.Item(12, 1).Value = "Recurring meeting"
$ws.Cells.Item(12, 2).Value = "May. 7"
$ws.Cells.Item(12, 3).Value = 1

# Update the selection to match the post-edit state
$ws.Range("D13").Select()
